$wb = $excel.ActiveWorkbook

# --- Charge Station AvailabilityZone sheet: insert a new "Abundance power" column ---
$wsZone = $wb.Worksheets.Item("Charge Station AvailabilityZone")

# Insert a new column before the existing "Location" column (column G / index 7),
# shifting Location / Zone priority / Station priority one column to the right.
$wsZone.Columns.Item(7).Insert()

# Populate the new column's header and data.
$wsZone.Cells.Item(1, 7).Value = "Abundance power"
$wsZone.Cells.Item(2, 7).Value = 60
$wsZone.Cells.Item(3, 7).Value = 10
$wsZone.Cells.Item(4, 7).Value = 90
$wsZone.Cells.Item(5, 7).Value = 100

# Give the new column its own width.
$wsZone.Columns.Item(7).ColumnWidth = 16.6

# Make this sheet the active tab, with G2 selected (matching the author's saved view).
# (This also clears the "EV Counts" sheet's previous tabSelected flag.)
$wsZone.Activate()
$wsZone.Range("G2").Select()
